$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1464.7273
$ws.Range("I28").Value = 1464.7273
$ws.Range("K28").Value = 1464.7273
$ws.Range("M28").Value = -979.7273
$ws.Range("H40").Value = 2197.077
$ws.Range("I40").Value = 2131.1875
$ws.Range("J40").Value = 2302.5
$ws.Range("K40").Value = 2131.1875
$ws.Range("L40").Value = 2302.5
$ws.Range("M40").Value = -1956.1875
$ws.Range("N40").Value = -2652.5
$ws.Range("H74").Value = 17708.834
$ws.Range("I74").Value = 25907.3
$ws.Range("J74").Value = 11852.786
$ws.Range("K74").Value = 25907.3
$ws.Range("L74").Value = 11852.786
$ws.Range("M74").Value = -24971.3
$ws.Range("N74").Value = -13724.786
$ws.Range("H77").Value = 17708.834
$ws.Range("I77").Value = 25907.3
$ws.Range("J77").Value = 11852.786
$ws.Range("K77").Value = 129536.5
$ws.Range("L77").Value = 59263.93
$ws.Range("M77").Value = -124856.5
$ws.Range("N77").Value = -68623.92999999999
$ws.Range("H107").Value = 522
$ws.Range("I107").Value = 502.5
$ws.Range("K107").Value = 502.5
$ws.Range("M107").Value = 1417.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1152.875
$ws.Range("I2").Value = 603.2857
$ws.Range("K2").Value = 603.2857
$ws.Range("M2").Value = -490.2857
$ws.Range("H45").Value = 1418.6154
$ws.Range("I45").Value = 1184.3
$ws.Range("K45").Value = 1184.3
$ws.Range("M45").Value = -807.3
$ws.Range("H61").Value = 3254.6538
$ws.Range("I61").Value = 1506.7368
$ws.Range("K61").Value = 1506.7368
$ws.Range("M61").Value = -1294.7368
$ws.Range("H102").Value = 875.2143
$ws.Range("I102").Value = 888.8461
$ws.Range("J102").Value = 698
$ws.Range("K102").Value = 888.8461
$ws.Range("L102").Value = 698
$ws.Range("M102").Value = 733.1539
$ws.Range("N102").Value = -3942
$ws.Range("H116").Value = 1152.875
$ws.Range("I116").Value = 603.2857
$ws.Range("K116").Value = 603.2857
$ws.Range("M116").Value = 1690.7143
$ws.Range("H122").Value = 2016.5
$ws.Range("I122").Value = 1566.6
$ws.Range("K122").Value = 4699.799999999999
$ws.Range("M122").Value = -2249.799999999999
$ws.Range("H132").Value = 1518.375
$ws.Range("I132").Value = 1449.7858
$ws.Range("K132").Value = 4349.357400000001
$ws.Range("M132").Value = -1819.357400000001
$ws.Range("H135").Value = 550000
$ws.Range("J135").Value = 550000
$ws.Range("L135").Value = 550000
$ws.Range("N135").Value = -560140
$ws.Range("H136").Value = 3254.6538
$ws.Range("I136").Value = 1506.7368
$ws.Range("K136").Value = 4520.2104
$ws.Range("M136").Value = -1970.2104

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1152.875
$ws.Range("I3").Value = 603.2857
$ws.Range("K3").Value = 603.2857
$ws.Range("M3").Value = -489.2857
$ws.Range("H99").Value = 2830.7646
$ws.Range("I99").Value = 1786.3
$ws.Range("J99").Value = 4322.857
$ws.Range("K99").Value = 1786.3
$ws.Range("L99").Value = 4322.857
$ws.Range("M99").Value = -288.3
$ws.Range("N99").Value = -7318.857

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1126.3636
$ws.Range("I16").Value = 1111.2632
$ws.Range("J16").Value = 1222
$ws.Range("K16").Value = 1111.2632
$ws.Range("L16").Value = 1222
$ws.Range("M16").Value = -824.2632000000001
$ws.Range("N16").Value = -1796
$ws.Range("H113").Value = 1126.3636
$ws.Range("I113").Value = 1111.2632
$ws.Range("J113").Value = 1222
$ws.Range("K113").Value = 1111.2632
$ws.Range("L113").Value = 1222
$ws.Range("M113").Value = 1058.7368
$ws.Range("N113").Value = -5562
$ws.Range("H134").Value = 2693
$ws.Range("I134").Value = 2528.2144
$ws.Range("K134").Value = 7584.6432
$ws.Range("M134").Value = -5049.6432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4359
$ws.Range("I70").Value = 4292.3335
$ws.Range("K70").Value = 4292.3335
$ws.Range("M70").Value = -4022.3335
$ws.Range("H73").Value = 4359
$ws.Range("I73").Value = 4292.3335
$ws.Range("K73").Value = 4292.3335
$ws.Range("M73").Value = -3356.3335
$ws.Range("H80").Value = 3022
$ws.Range("I80").Value = 1287.25
$ws.Range("K80").Value = 1287.25
$ws.Range("M80").Value = -289.25
$ws.Range("H83").Value = 3022
$ws.Range("I83").Value = 1287.25
$ws.Range("K83").Value = 6436.25
$ws.Range("M83").Value = -1444.25
$ws.Range("H122").Value = 2350
$ws.Range("I122").Value = 1303.6666
$ws.Range("K122").Value = 3910.9998
$ws.Range("M122").Value = -1460.9998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2879.5454
$ws.Range("I46").Value = 2493.25
$ws.Range("K46").Value = 2493.25
$ws.Range("M46").Value = -2305.25
$ws.Range("H61").Value = 4045.4119
$ws.Range("I61").Value = 4187.273
$ws.Range("J61").Value = 3785.3333
$ws.Range("K61").Value = 4187.273
$ws.Range("L61").Value = 3785.3333
$ws.Range("M61").Value = -3985.273
$ws.Range("N61").Value = -4189.3333
$ws.Range("H82").Value = 3199.6
$ws.Range("I82").Value = 2666.3333
$ws.Range("J82").Value = 3999.5
$ws.Range("K82").Value = 2666.3333
$ws.Range("L82").Value = 3999.5
$ws.Range("M82").Value = -2305.3333
$ws.Range("N82").Value = -4721.5
$ws.Range("H85").Value = 3199.6
$ws.Range("I85").Value = 2666.3333
$ws.Range("J85").Value = 3999.5
$ws.Range("K85").Value = 2666.3333
$ws.Range("L85").Value = 3999.5
$ws.Range("M85").Value = -1418.3333
$ws.Range("N85").Value = -6495.5
$ws.Range("H113").Value = 4045.4119
$ws.Range("I113").Value = 4187.273
$ws.Range("J113").Value = 3785.3333
$ws.Range("K113").Value = 4187.273
$ws.Range("L113").Value = 3785.3333
$ws.Range("M113").Value = -2017.273
$ws.Range("N113").Value = -8125.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 10000533
$ws.Range("I100").Value = 25000394
$ws.Range("J100").Value = 625
$ws.Range("K100").Value = 50000788
$ws.Range("L100").Value = 1250
$ws.Range("M100").Value = -50000247
$ws.Range("N100").Value = -2332
$ws.Range("H107").Value = 1345.8889
$ws.Range("I107").Value = 519.5
$ws.Range("J107").Value = 2998.6667
$ws.Range("K107").Value = 1558.5
$ws.Range("L107").Value = 8996.000100000001
$ws.Range("M107").Value = 361.5
$ws.Range("N107").Value = -12836.0001
$ws.Range("H122").Value = 761.06665
$ws.Range("I122").Value = 769.4286
$ws.Range("K122").Value = 2308.2858
$ws.Range("M122").Value = 141.7142000000003
